$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & volume change) per row.
# For numeric-looking Price strings we force a Text format before writing
# the value (so Excel keeps the exact original text, e.g. trailing zeros,
# instead of re-interpreting it as a number), then restore the default
# "Normal" style so no stray number-format style is left on the cell.

# Row 2
$ws.Range("D2").Value = "61.308.53"
$ws.Range("E2").Value = "  -2.44%  "

# Row 3
$ws.Range("D3").Value = "3.391.03"
$ws.Range("E3").Value = "  -2.24%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "403.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.90%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.13%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.679"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.08%  "

# Row 10
$ws.Range("E10").Value = "  -7.15%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.84%  "

# Row 12
$ws.Range("E12").Value = "  -1.15%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.31%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.13%  "

# Row 15
$ws.Range("D15").Value = "3.399.83"
$ws.Range("E15").Value = "  -1.86%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.69"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.32%  "

# Row 17
$ws.Range("D17").Value = "61.358.94"
$ws.Range("E17").Value = "  -2.50%  "

# Row 18
$ws.Range("E18").Value = "  -3.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000140"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.76%  "

# Row 20
$ws.Range("E20").Value = "  -6.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "82.90"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "310.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.87%  "

# Row 24
$ws.Range("E24").Value = "  -1.33%  "

# Row 25
$ws.Range("E25").Value = "  +9.96%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.63%  "

# Row 28
$ws.Range("E28").Value = "  -2.71%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.31%  "

# Row 30
$ws.Range("E30").Value = "  -3.89%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "43.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "

# Row 32
$ws.Range("E32").Value = "  -3.31%  "

# Row 33
$ws.Range("E33").Value = "  +0.11%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.57%  "

# Row 35
$ws.Range("E35").Value = "  -3.11%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "51.27"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.59%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.08%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.15%  "

# Row 39
$ws.Range("E39").Value = "  -3.27%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.317"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "139.39"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "

# Row 42
$ws.Range("E42").Value = "  -1.70%  "

# Row 43
$ws.Range("E43").Value = "  -2.38%  "

# Row 44
$ws.Range("E44").Value = "  -1.73%  "

# Row 45
$ws.Range("E45").Value = "  -6.28%  "

# Row 46
$ws.Range("E46").Value = "  -1.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.23%  "

# Row 48
$ws.Range("D48").Value = "2.088.80"
$ws.Range("E48").Value = "  -4.30%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.88%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +20.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.30%  "
